$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 33 with the next record (Mac-Address / Document Types entry)
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"
$ws.Range("G33").Value = "now()"

# Update the selected cell to match the new view state
$ws.Range("D26").Select()
